$wb = $excel.ActiveWorkbook

# --- About sheet: update the "last updated" date cell (C1) ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45387

# --- BAU Emissions sheet: rename the " : NoSettings" suffix to " : test" ---
# These row-label strings live in column A and are unique to this sheet, so a
# plain find/replace over the used rows reproduces the shared-string edit.
$wsBau = $wb.Worksheets.Item("BAU Emissions")
$lastRow = $wsBau.UsedRange.Rows.Count
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $wsBau.Cells.Item($r, 1)
    $txt = $cell.Value2
    if ($txt -ne $null -and $txt -like "*: NoSettings") {
        $cell.Value = ($txt -replace ': NoSettings$', ': test')
    }
}

# --- BAU Emissions sheet: refresh the "natural gas if,iron and steel" row (row 94) ---
$wsBau.Range("M94").Value = 1001080
$wsBau.Range("N94").Value = 2002150
$wsBau.Range("O94").Value = 3003230
$wsBau.Range("P94").Value = 4004300
$wsBau.Range("Q94:AE94").Value = 5005380

# --- View state: BAU Emissions selection moves, then focus returns to About ---
$wsBau.Activate()
$wsBau.Range("A30:AE280").Select()

$wsAbout.Activate()
